# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 290
$wsExhibit.Range("F4").Value  = 3534
$wsExhibit.Range("F5").Value  = 2199
$wsExhibit.Range("F7").Value  = 170
$wsExhibit.Range("F8").Value  = 69
$wsExhibit.Range("F9").Value  = 59
$wsExhibit.Range("F10").Value = 1304
$wsExhibit.Range("F12").Value = 1792
$wsExhibit.Range("F13").Value = 133

# --- Sheet "全部类型" (index 4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 290
$wsAll.Range("F4").Value  = 3534
$wsAll.Range("F5").Value  = 2199
$wsAll.Range("F8").Value  = 170
$wsAll.Range("F9").Value  = 69
$wsAll.Range("F10").Value = 59
$wsAll.Range("F13").Value = 1304
$wsAll.Range("F15").Value = 1792
$wsAll.Range("F16").Value = 133
